# FA60_TestData_PostMassAdditionsSubmission_21C.xlsx
# "Add files via upload" / "Anu - FA files uploaded"
#
# The sheet's Input_Value tab carried a live URL + credential pair
# (L2/M2/N2 -> https://edrx.fa.us2.oraclecloud.com, IBM_IMPLEMENTATION_USER,
# Oracle1234) with a hyperlink on the URL cell. This re-upload scrubs those
# secrets: the three cells are cleared and the hyperlink is removed, which
# also drops the now-unused shared strings.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Input_Value")

# Drop the hyperlink that was attached to the URL cell (L2).
$ws1.Hyperlinks.Delete()

# Scrub the credentials that used to live in L2:N2 (URL / UserName / Password).
$ws1.Range("L2").ClearContents()
$ws1.Range("M2").ClearContents()
$ws1.Range("N2").ClearContents()

# Leave the selection where the edit happened.
$ws1.Range("L2:N2").Select()
